# Add solutions for LeetCode problems 1133, 1150 and 2197 to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 162: 1133. Largest Unique Number ---------------------------------
$ws.Range("A161:I161").Copy()
$ws.Range("A162:I162").PasteSpecial(-4122)  # xlPasteFormats (keeps existing style ids, avoids new numFmts)
$ws.Rows.Item(162).RowHeight = 34

$ws.Range("A162").Value = 1133
$ws.Range("B162").Value = "Largest Unique Number"
$ws.Range("C162").Value = "#array  #hash-table "
$ws.Range("D162").Value = "easy"
$ws.Range("E162").Value = 1
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 45915
$ws.Range("I162").Value = 45915

# --- Row 163: 1150. Check If a Number Is Majority Element in a Sorted Array
$ws.Range("A161:I161").Copy()
$ws.Range("A163:I163").PasteSpecial(-4122)
$ws.Rows.Item(163).RowHeight = 51

$ws.Range("A163").Value = 1150
$ws.Range("B163").Value = "Check If a Number Is Majority Element in a Sorted Array"
$ws.Range("C163").Value = "#array #binary-search "
$ws.Range("D163").Value = "easy"
$ws.Range("E163").Value = 1
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 45915
$ws.Range("I163").Value = 45915

# --- Row 164: 2197. Replace Non-Coprime Numbers in Array -------------------
$ws.Range("A161:I161").Copy()
$ws.Range("A164:I164").PasteSpecial(-4122)
$ws.Rows.Item(164).RowHeight = 34

$ws.Range("A164").Value = 2197
$ws.Range("B164").Value = "Replace Non-Coprime Numbers in Array"
$ws.Range("C164").Value = "#array #greedy #stack "
$ws.Range("D164").Value = "hard"
$ws.Range("E164").Value = 0
$ws.Range("F164").Value = 1
$ws.Range("G164").Value = 45
$ws.Range("H164").Value = 45916
$ws.Range("I164").Value = 45916

# --- View state: selection moves to H167 after the new rows ---------------
$ws.Range("H167").Select()
